# Generate Report for Handoff
#
# The "a.md" file has now been handed off for localization ("Ready for
# handoff"), while the "b.md" file (previously "In Translation") moves up
# to take the first data row. In effect the two data rows swap positions,
# and the row that now holds "a.md" gets a refreshed status/date and new
# handoff/error info.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Row 2 becomes the "b.md" row (previously row 3)
$ws.Range("A2").Value = "b.md"
$ws.Range("B2").Value = "e2e\b.md"
$ws.Range("C2").Value = ".md"
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = "In Translation"
$ws.Range("F2").Value = "In Translation"
$ws.Range("G2").Value = "2017-01-03 07:15:46"

# Row 3 becomes the "a.md" row, now ready for handoff
$ws.Range("A3").Value = "a.md"
$ws.Range("B3").Value = "e2e\a.md"
$ws.Range("C3").Value = ".md"
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = "Ready for handoff"
$ws.Range("F3").Value = "Ready for handoff"
$ws.Range("G3").Value = "2017-01-03 07:16:48"

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = "b.md"
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "In Translation"
$ws.Range("D2").Value = "e2e"
$ws.Range("E2").Value = "ht"
$ws.Range("F2").Value = "'False"
$ws.Range("G2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$ws.Range("H2").Value = "2017-01-03 07:16:37"
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = "a.md"
$ws.Range("K2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$ws.Range("L2").Value = "2017-01-03 07:14:55"
$ws.Range("M2").Value = ""
$ws.Range("N2").Value = ""
$ws.Range("O2").Value = "'True"
$ws.Range("P2").Value = ""
$ws.Range("Q2").Value = "'False"
$ws.Range("R2").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test1/blob/a2a32c0a5631a13868300f9802c4f19155acb3b9/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test1/blob/a4b80ca38a28cdfb7ac57ae17c50f2577c2c14dc/e2e/b.md."

$ws.Range("A3").Value = "a.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "e2e"
$ws.Range("E3").Value = "ht"
$ws.Range("F3").Value = "'False"
$ws.Range("G3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$ws.Range("H3").Value = "2017-01-03 07:16:37"
$ws.Range("I3").Value = ""
$ws.Range("J3").Value = "a.md"
$ws.Range("K3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$ws.Range("L3").Value = "2017-01-03 07:14:55"
$ws.Range("M3").Value = ""
$ws.Range("N3").Value = ""
$ws.Range("O3").Value = "'True"
$ws.Range("P3").Value = ""
$ws.Range("Q3").Value = "'False"
$ws.Range("R3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test1/blob/a2a32c0a5631a13868300f9802c4f19155acb3b9/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test1/blob/e869d653495c60f9140b6df56e719a6744700ae8/e2e/a.md."

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = "b.md"
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "In Translation"
$ws.Range("D2").Value = "e2e"
$ws.Range("E2").Value = "ht"
$ws.Range("F2").Value = "'False"
$ws.Range("G2").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$ws.Range("H2").Value = "2017-01-03 07:15:46"
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = "a.md"
$ws.Range("K2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$ws.Range("L2").Value = "2017-01-03 07:15:07"
$ws.Range("M2").Value = ""
$ws.Range("N2").Value = ""
$ws.Range("O2").Value = "'True"
$ws.Range("P2").Value = ""
$ws.Range("Q2").Value = "'False"
$ws.Range("R2").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test1/blob/a2a32c0a5631a13868300f9802c4f19155acb3b9/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test1/blob/a4b80ca38a28cdfb7ac57ae17c50f2577c2c14dc/e2e/b.md."

$ws.Range("A3").Value = "a.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "e2e"
$ws.Range("E3").Value = "ht"
$ws.Range("F3").Value = "'False"
$ws.Range("G3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$ws.Range("H3").Value = "2017-01-03 07:16:48"
$ws.Range("I3").Value = ""
$ws.Range("J3").Value = "a.md"
$ws.Range("K3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$ws.Range("L3").Value = "2017-01-03 07:15:07"
$ws.Range("M3").Value = ""
$ws.Range("N3").Value = ""
$ws.Range("O3").Value = "'True"
$ws.Range("P3").Value = ""
$ws.Range("Q3").Value = "'False"
$ws.Range("R3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test1/blob/a2a32c0a5631a13868300f9802c4f19155acb3b9/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test1/blob/e869d653495c60f9140b6df56e719a6744700ae8/e2e/a.md."

# ---------------------------------------------------------------------
# Column widths: narrow the long-text columns (Overview E/F, zh-cn & de-de C)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E1:F1").ColumnWidth = 17.2159881591797

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C1").ColumnWidth = 17.2159881591797

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C1").ColumnWidth = 17.2159881591797
